$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 values
$ws.Range("E4").Value = 27714428
$ws.Range("F4").Value = 50.34
$ws.Range("G4").Value = 4.14

# Row 5 values
$ws.Range("E5").Value = 33714428
$ws.Range("F5").Value = 33714428
$ws.Range("G5").Value = 1203

# Apply new number format (#,##0.000) to E4:E5 (per target diff, cellXf index 6 ends up with this numFmt)
$ws.Range("E4:E5").NumberFormat = "#,##0.000"
